# Auto-generated Excel COM-interop script applying the Aegis_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 125950
$ws.Range("J97").Value = 125950
$ws.Range("L97").Value = 377850
$ws.Range("N97").Value = -378842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1734.8948
$ws.Range("J112").Value = 1841.4375
$ws.Range("L112").Value = 5524.3125
$ws.Range("N112").Value = -7740.3125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3082.5715
$ws.Range("I116").Value = 2599.6667
$ws.Range("J116").Value = 3444.75
$ws.Range("K116").Value = 2599.6667
$ws.Range("L116").Value = 3444.75
$ws.Range("M116").Value = 842.3332999999998
$ws.Range("N116").Value = -10328.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1591.55
$ws.Range("I137").Value = 1607.7059
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 4823.1177
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -2273.1177
$ws.Range("N137").Value = -9600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26408.54
$ws.Range("I32").Value = 5886.147
$ws.Range("J32").Value = 200848.88
$ws.Range("K32").Value = 5886.147
$ws.Range("L32").Value = 200848.88
$ws.Range("M32").Value = -5599.147
$ws.Range("N32").Value = -201422.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1981.9131
$ws.Range("I61").Value = 1839.3334
$ws.Range("J61").Value = 2249.25
$ws.Range("K61").Value = 1839.3334
$ws.Range("L61").Value = 2249.25
$ws.Range("M61").Value = -1627.3334
$ws.Range("N61").Value = -2673.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2955.075
$ws.Range("I132").Value = 2833.9033
$ws.Range("J132").Value = 3372.4443
$ws.Range("K132").Value = 8501.7099
$ws.Range("L132").Value = 10117.3329
$ws.Range("M132").Value = -5971.7099
$ws.Range("N132").Value = -15177.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1981.9131
$ws.Range("I136").Value = 1839.3334
$ws.Range("J136").Value = 2249.25
$ws.Range("K136").Value = 5518.0002
$ws.Range("L136").Value = 6747.75
$ws.Range("M136").Value = -2968.0002
$ws.Range("N136").Value = -11847.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2610.5264
$ws.Range("I6").Value = 2800
$ws.Range("J6").Value = 2200
$ws.Range("K6").Value = 2800
$ws.Range("L6").Value = 2200
$ws.Range("M6").Value = -2687
$ws.Range("N6").Value = -2426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.18182
$ws.Range("I7").Value = 140
$ws.Range("J7").Value = 275.25
$ws.Range("K7").Value = 140
$ws.Range("L7").Value = 275.25
$ws.Range("M7").Value = -27
$ws.Range("N7").Value = -501.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 15000
$ws.Range("I17").Value = 15000
$ws.Range("K17").Value = 15000
$ws.Range("M17").Value = -14826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22154
$ws.Range("I31").Value = 1470.3478
$ws.Range("J31").Value = 32275.787
$ws.Range("K31").Value = 1470.3478
$ws.Range("L31").Value = 32275.787
$ws.Range("M31").Value = -1175.3478
$ws.Range("N31").Value = -32865.787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 22154
$ws.Range("I34").Value = 1470.3478
$ws.Range("J34").Value = 32275.787
$ws.Range("K34").Value = 1470.3478
$ws.Range("L34").Value = 32275.787
$ws.Range("M34").Value = -1268.3478
$ws.Range("N34").Value = -32679.787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9987.286
$ws.Range("I41").Value = 3417
$ws.Range("J41").Value = 14915
$ws.Range("K41").Value = 3417
$ws.Range("L41").Value = 14915
$ws.Range("M41").Value = -2989
$ws.Range("N41").Value = -15771

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 16792
$ws.Range("J50").Value = 16792
$ws.Range("L50").Value = 16792
$ws.Range("N50").Value = -18042

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8024.5
$ws.Range("J51").Value = 8024.5
$ws.Range("L51").Value = 8024.5
$ws.Range("N51").Value = -9496.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2327.923
$ws.Range("I58").Value = 1993.75
$ws.Range("J58").Value = 2862.6
$ws.Range("K58").Value = 1993.75
$ws.Range("L58").Value = 2862.6
$ws.Range("M58").Value = -1790.75
$ws.Range("N58").Value = -3268.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 24796
$ws.Range("J59").Value = 24796
$ws.Range("L59").Value = 24796
$ws.Range("N59").Value = -27086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 17010.666
$ws.Range("J60").Value = 19490
$ws.Range("L60").Value = 19490
$ws.Range("N60").Value = -20512

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 8024.5
$ws.Range("J61").Value = 8024.5
$ws.Range("L61").Value = 8024.5
$ws.Range("N61").Value = -8720.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 18241.809
$ws.Range("J68").Value = 18241.809
$ws.Range("L68").Value = 18241.809
$ws.Range("N68").Value = -19739.809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 18241.809
$ws.Range("J71").Value = 18241.809
$ws.Range("L71").Value = 54725.427
$ws.Range("N71").Value = -62213.427

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 29500
$ws.Range("J74").Value = 29500
$ws.Range("L74").Value = 29500
$ws.Range("N74").Value = -31248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 29500
$ws.Range("J77").Value = 29500
$ws.Range("L77").Value = 88500
$ws.Range("N77").Value = -97236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9513.77
$ws.Range("I99").Value = 2293
$ws.Range("J99").Value = 12723
$ws.Range("K99").Value = 2293
$ws.Range("L99").Value = 12723
$ws.Range("M99").Value = -795
$ws.Range("N99").Value = -15719

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 164000
$ws.Range("J121").Value = 164000
$ws.Range("L121").Value = 164000
$ws.Range("N121").Value = -166620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9513.77
$ws.Range("I126").Value = 2293
$ws.Range("J126").Value = 12723
$ws.Range("K126").Value = 6879
$ws.Range("L126").Value = 38169
$ws.Range("M126").Value = -4409
$ws.Range("N126").Value = -43109

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2327.923
$ws.Range("I136").Value = 1993.75
$ws.Range("J136").Value = 2862.6
$ws.Range("K136").Value = 5981.25
$ws.Range("L136").Value = 8587.799999999999
$ws.Range("M136").Value = -3431.25
$ws.Range("N136").Value = -13687.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 770.59
$ws.Range("I131").Value = 360.2857
$ws.Range("J131").Value = 801.4731399999999
$ws.Range("K131").Value = 1080.8571
$ws.Range("L131").Value = 2404.41942
$ws.Range("M131").Value = 3959.1429
$ws.Range("N131").Value = -12484.41942

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 8586.5
$ws.Range("I34").Value = 9000
$ws.Range("J34").Value = 8173
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 8173
$ws.Range("M34").Value = -8732
$ws.Range("N34").Value = -8709

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 8586.5
$ws.Range("I76").Value = 9000
$ws.Range("J76").Value = 8173
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 8173
$ws.Range("M76").Value = -8685
$ws.Range("N76").Value = -8803

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 8586.5
$ws.Range("I79").Value = 9000
$ws.Range("J79").Value = 8173
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 8173
$ws.Range("M79").Value = -7908
$ws.Range("N79").Value = -10357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3189.111
$ws.Range("I99").Value = 3189.111
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3189.111
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -943.1109999999999
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 778.6667
$ws.Range("I122").Value = 762.8570999999999
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2288.5713
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 161.4287000000004
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5459.4
$ws.Range("I7").Value = 5477.5557
$ws.Range("J7").Value = 5432.1665
$ws.Range("K7").Value = 5477.5557
$ws.Range("L7").Value = 5432.1665
$ws.Range("M7").Value = -5365.5557
$ws.Range("N7").Value = -5656.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 40113.69
$ws.Range("I40").Value = 101085.4
$ws.Range("J40").Value = 2006.375
$ws.Range("K40").Value = 101085.4
$ws.Range("L40").Value = 2006.375
$ws.Range("M40").Value = -100949.4
$ws.Range("N40").Value = -2278.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 487078.94
$ws.Range("I46").Value = 304.7
$ws.Range("J46").Value = 929601
$ws.Range("K46").Value = 304.7
$ws.Range("L46").Value = 929601
$ws.Range("M46").Value = -116.7
$ws.Range("N46").Value = -929977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3016.25
$ws.Range("J122").Value = 3039
$ws.Range("L122").Value = 9117
$ws.Range("N122").Value = -14017

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5459.4
$ws.Range("I126").Value = 5477.5557
$ws.Range("J126").Value = 5432.1665
$ws.Range("K126").Value = 16432.6671
$ws.Range("L126").Value = 16296.4995
$ws.Range("M126").Value = -13962.6671
$ws.Range("N126").Value = -21236.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 14533.333
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 14533.333
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 14533.333
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -15359.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 27208.166
$ws.Range("J42").Value = 27208.166
$ws.Range("L42").Value = 27208.166
$ws.Range("N42").Value = -27964.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 2500
$ws.Range("I43").Value = 2500
$ws.Range("K43").Value = 2500
$ws.Range("M43").Value = -2351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 36553.668
$ws.Range("J118").Value = 36553.668
$ws.Range("L118").Value = 36553.668
$ws.Range("N118").Value = -39867.668
